$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4-5, shifting existing data (rows 4-29) down to rows 6-31
$ws.Rows("4:5").Insert()

# Match formatting of index column (A) for the two new rows to the existing style
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 4: "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.006014271474516
$ws.Range("D4").Value = 0.8845528534085375
$ws.Range("E4").Value = 1.041819536877271
$ws.Range("F4").Value = 0.8845528534085375
$ws.Range("G4").Value = 1.031380684084285
$ws.Range("H4").Value = 1.006014271474516
$ws.Range("I4").Value = 1.121135329344442
$ws.Range("J4").Value = 0.9271208647971515
$ws.Range("K4").Value = 1.006014271474516
$ws.Range("L4").Value = 1.041819536877271
$ws.Range("M4").Value = 0.9631861951429044
$ws.Range("N4").Value = 0.9631861951429044
$ws.Range("O4").Value = 0.9511644183609868
$ws.Range("P4").Value = 0.9774622205867748
$ws.Range("Q4").Value = 0.9774622205867747
$ws.Range("R4").Value = 0.9846002333087099
$ws.Range("S4").Value = 0.9846002333087099
$ws.Range("T4").Value = 1.002003923331034

# New row 5: "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.008745944250528
$ws.Range("D5").Value = 0.9171184436391582
$ws.Range("E5").Value = 1.027508283472786
$ws.Range("F5").Value = 0.9171184436391582
$ws.Range("G5").Value = 1.022038211484679
$ws.Range("H5").Value = 1.008745944250528
$ws.Range("I5").Value = 1.078010048043325
$ws.Range("J5").Value = 0.9495934548891013
$ws.Range("K5").Value = 1.008745944250528
$ws.Range("L5").Value = 1.027508283472786
$ws.Range("M5").Value = 0.9723133635559721
$ws.Range("N5").Value = 0.9723133635559721
$ws.Range("O5").Value = 0.9647400606670152
$ws.Range("P5").Value = 0.9844575571208241
$ws.Range("Q5").Value = 0.9844575571208241
$ws.Range("R5").Value = 0.9905296539032501
$ws.Range("S5").Value = 0.9905296539032501
$ws.Range("T5").Value = 1.00050239762993

# Rename "Thomas Hex" -> "Matthies Hex" (row shifted from 9 to 11 after insert)
$ws.Range("B11").Value = "Matthies Hex"
